# Weekly fruit/vegetable price update:
# Two new daily records (fecha serial 44587) are inserted right before the
# existing row 116, pushing all subsequent rows down by two positions
# (former row 116 -> 118, ... former row 214 -> 216).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 116-117; everything below (old 116..214) shifts
# down to 118..216, carrying its formatting (incl. the date style on column D).
$ws.Rows("116:117").Insert()

# --- New row 116: Americana (o) / Primera ---
$ws.Range("A116").Value = 2
$ws.Range("B116").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C116").Value = "Coquimbo"
$ws.Range("D116").Value = 44587
$ws.Range("E116").Value = 4
$ws.Range("F116").Value = 100112021
$ws.Range("G116").Value = "Ají"
$ws.Range("H116").Value = "Americana (o)"
$ws.Range("I116").Value = "Primera"
$ws.Range("J116").Value = 500
$ws.Range("K116").Value = 8000
$ws.Range("L116").Value = 10000
$ws.Range("M116").Value = 9000
$ws.Range("N116").Value = "$/caja 25 kilos"
$ws.Range("O116").Value = "Provincia de Limarí"
$ws.Range("P116").Value = 360
$ws.Range("Q116").Value = 25
$ws.Range("R116").Value = "Hortaliza"

# --- New row 117: Inferno / Primera ---
$ws.Range("A117").Value = 2
$ws.Range("B117").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C117").Value = "Coquimbo"
$ws.Range("D117").Value = 44587
$ws.Range("E117").Value = 4
$ws.Range("F117").Value = 100112021
$ws.Range("G117").Value = "Ají"
$ws.Range("H117").Value = "Inferno"
$ws.Range("I117").Value = "Primera"
$ws.Range("J117").Value = 300
$ws.Range("K117").Value = 12000
$ws.Range("L117").Value = 13000
$ws.Range("M117").Value = 12500
$ws.Range("N117").Value = "$/caja 25 kilos"
$ws.Range("O117").Value = "Provincia de Limarí"
$ws.Range("P117").Value = 500
$ws.Range("Q117").Value = 25
$ws.Range("R117").Value = "Hortaliza"
